# Weekly update for "Hortaliza, Vega Modelo de Temuco - Coliflor":
# a new daily-price record is inserted as row 739 (pushing the existing
# rows 739-789 down to 740-790), extending the used range to A1:R790.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at 739, shifting rows 739..789 down to 740..790.
$ws.Rows.Item(739).Insert()

# Populate the new row 739 with the latest observation.
$ws.Cells.Item(739, 1).Value2  = 10
$ws.Cells.Item(739, 2).Value2  = "Vega Modelo de Temuco"
$ws.Cells.Item(739, 3).Value2  = "La Araucanía"
$ws.Cells.Item(739, 4).Value2  = 45265
$ws.Cells.Item(739, 5).Value2  = 9
$ws.Cells.Item(739, 6).Value2  = 100112008
$ws.Cells.Item(739, 7).Value2  = "Coliflor"
$ws.Cells.Item(739, 8).Value2  = "Sin especificar"
$ws.Cells.Item(739, 9).Value2  = "Primera"
$ws.Cells.Item(739, 10).Value2 = 1200
$ws.Cells.Item(739, 11).Value2 = 1400
$ws.Cells.Item(739, 12).Value2 = 1400
$ws.Cells.Item(739, 13).Value2 = 1400
$ws.Cells.Item(739, 14).Value2 = "`$/unidad"
$ws.Cells.Item(739, 15).Value2 = "Región Metropolitana"
$ws.Cells.Item(739, 16).Value2 = 1400
$ws.Cells.Item(739, 17).Value2 = 1
$ws.Cells.Item(739, 18).Value2 = "Hortaliza"
